$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Through 2022-08-17"

# Update header text in I1
$ws.Range("I1").Value = "2022 (through 08-17)"

# Update data values
$ws.Range("I9").Value = 92
$ws.Range("I14").Value = 1063
